$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("CityResaleNum")

$row = 62

# Columns A and D contain values that look like dates/numbers to Excel's
# auto-detection (e.g. "2025-02-10" and "06"), so force them to be treated
# as plain text the same way the rest of the sheet stores them, without
# leaving a residual number format applied to the cell.
$cellA = $ws.Cells.Item($row, 1)
$cellA.NumberFormat = "@"
$cellA.Value = "2025-02-10"
$cellA.ClearFormats()

$ws.Cells.Item($row, 2).Value = "22:45:00"
$ws.Cells.Item($row, 3).Value = "Monday"

$cellD = $ws.Cells.Item($row, 4)
$cellD.NumberFormat = "@"
$cellD.Value = "06"
$cellD.ClearFormats()

$ws.Cells.Item($row, 5).Value = 127600
$ws.Cells.Item($row, 6).Value = 141935
$ws.Cells.Item($row, 7).Value = 169314
$ws.Cells.Item($row, 8).Value = 158449
$ws.Cells.Item($row, 9).Value = -1
$ws.Cells.Item($row, 10).Value = 144289
$ws.Cells.Item($row, 11).Value = -1
$ws.Cells.Item($row, 12).Value = -1
$ws.Cells.Item($row, 13).Value = 191760
$ws.Cells.Item($row, 14).Value = 115056
$ws.Cells.Item($row, 15).Value = 44926
$ws.Cells.Item($row, 16).Value = 28495
$ws.Cells.Item($row, 17).Value = 64815
$ws.Cells.Item($row, 18).Value = -1
$ws.Cells.Item($row, 19).Value = 43416
$ws.Cells.Item($row, 20).Value = -1
